# Apply updated renta_sobre_pbi data to both sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet 1: metodo_directo_pbi ---
$ws1 = $wb.Worksheets.Item("metodo_directo_pbi")

$sheet1Data = @(
    @(1998, 0.01765449455060211, 0.00860926459579257),
    @(1999, 0.02304916703296808, 0.01088031701055575),
    @(2000, 0.04177064086469819, 0.0202858319208292),
    @(2001, 0.03105534818247203, 0.01514600183863189),
    @(2001, 0.03088050444735142, 0.01506072881196833),
    @(2002, 0.01700036049260459, 0.008624283194065118),
    @(2002, 0.01705558267174119, 0.008652297406568121),
    @(2003, 0.02464732289044107, 0.01309983083349207),
    @(2003, 0.02481915318002814, 0.01319115708972145),
    @(2004, 0.04175519104138919, 0.0215913466864755),
    @(2004, 0.04206247384189608, 0.02175024069009673),
    @(2005, 0.07141438879948278, 0.03667871049851687),
    @(2005, 0.07145197800366925, 0.03669801646138236),
    @(2006, 0.08344896089961341, 0.04244456922932892),
    @(2006, 0.08373715399541222, 0.04259115262203002),
    @(2007, 0.08430482041274182, 0.04262650293190438),
    @(2007, 0.08470165314417698, 0.04282715090798902),
    @(2008, 0.172193226252312, 0.08714518399900682),
    @(2008, 0.1727063788126813, 0.08740488512236823),
    @(2009, 0.1041068073572681, 0.04552886620401946),
    @(2009, 0.1046804337975665, 0.04577972935230339),
    @(2010, 0.1279509298381581, 0.05955901973153042),
    @(2010, 0.1281206089705064, 0.05963800253223649),
    @(2011, 0.1399080615736059, 0.0644607268047141),
    @(2011, 0.1397549709704236, 0.06439019240207063),
    @(2012, 0.1629011849790095, 0.07031272952471411),
    @(2012, 0.1623000029088483, 0.07005324244793318),
    @(2013, 0.1200507325766218, 0.04998135405429974),
    @(2013, 0.1181457300236252, 0.04918823430374016),
    @(2014, 0.1121826707978376, 0.04926404588224025),
    @(2014, 0.1086627605569151, 0.04771830786069983),
    @(2015, 0.02604033125670879, 0.01120208958366055),
    @(2016, -0.01247550390540515, -0.005477154261997007),
    @(2017, 0.02386056756701393, 0.01038795396274794),
    @(2018, 0.04196129166117218, 0.01942539529351357)
)

$row = 2
foreach ($rec in $sheet1Data) {
    $ws1.Cells.Item($row, 1).Value = $rec[0]
    $ws1.Cells.Item($row, 2).Value = $rec[1]
    $ws1.Cells.Item($row, 3).Value = $rec[2]
    $row = $row + 1
}

# --- Sheet 2: metodo_indirecto_pbi ---
$ws2 = $wb.Worksheets.Item("metodo_indirecto_pbi")

$sheet2Data = @(
    @(89, 0.01917489908140569, 0.009350694199501243),
    @(90, 0.02849490501326884, 0.01345096763308569),
    @(91, 0.04633405520449613, 0.02250204537520497),
    @(92, 0.0336147868540467, 0.01639426550638285),
    @(93, 0.02554804725602533, 0.01296052484811537),
    @(94, 0.03650522048047757, 0.01940219694281891),
    @(95, 0.05778549441431757, 0.02988051574502902),
    @(96, 0.08388447499942674, 0.04308339573503894),
    @(97, 0.1142173818900187, 0.05809428326681602),
    @(98, 0.110851883741452, 0.05604932344530512),
    @(99, 0.2102472334224802, 0.1064039175097236),
    @(100, 0.1418345778131168, 0.06202829267634967),
    @(101, 0.1553416108018514, 0.07230892401155507),
    @(102, 0.1679872634890108, 0.07739783524010668),
    @(103, 0.1919197830020315, 0.08283797195460843),
    @(104, 0.145324614135775, 0.06050376233471762),
    @(105, 0.1382790759194888, 0.06072405561575192),
    @(106, 0.04903898868567586, 0.02109570492531774),
    @(107, 0.01100118785254267, 0.004829881293009537),
    @(108, 0.05487073081006221, 0.02388856107284098),
    @(109, 0.05979267963830075, 0.02768018789819984)
)

foreach ($rec in $sheet2Data) {
    $r = $rec[0]
    $ws2.Cells.Item($r, 2).Value = $rec[1]
    $ws2.Cells.Item($r, 3).Value = $rec[2]
}
